$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the whole Price column to Text first so numeric-looking strings
# (e.g. "1.003") are written verbatim instead of being parsed into numbers.
# The format is cleared again afterwards so cells keep their original
# (unstyled) look, matching the source data which stored these as plain text.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.309.32'
$ws.Range("E2").Value = '  +1.19%  '

$ws.Range("D3").Value = '1.808.58'
$ws.Range("E3").Value = '  +3.34%  '

$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = '338.03'
$ws.Range("E5").Value = '  +0.73%  '

$ws.Range("D6").Value = '0.9992'
$ws.Range("E6").Value = '  -0.23%  '

$ws.Range("D7").Value = '0.4658'
$ws.Range("E7").Value = '  +21.20%  '

$ws.Range("D8").Value = '0.3814'
$ws.Range("E8").Value = '  +12.51%  '

$ws.Range("D9").Value = '45.51'
$ws.Range("E9").Value = '  -0.56%  '

$ws.Range("D10").Value = '1.159'
$ws.Range("E10").Value = '  +4.10%  '

$ws.Range("D11").Value = '0.07640'
$ws.Range("E11").Value = '  +6.09%  '

$ws.Range("D12").Value = '22.52'
$ws.Range("E12").Value = '  +0.29%  '

$ws.Range("E13").Value = '  -0.37%  '

$ws.Range("D14").Value = '6.352'
$ws.Range("E14").Value = '  +3.37%  '

$ws.Range("D15").Value = '7.465'
$ws.Range("E15").Value = '  +5.28%  '

$ws.Range("D16").Value = '1.807.31'
$ws.Range("E16").Value = '  +3.10%  '

$ws.Range("D17").Value = '0.00001096'
$ws.Range("E17").Value = '  +3.69%  '

$ws.Range("D18").Value = '0.06718'
$ws.Range("E18").Value = '  +1.69%  '

$ws.Range("D19").Value = '81.90'
$ws.Range("E19").Value = '  +3.52%  '

$ws.Range("D20").Value = '0.9992'
$ws.Range("E20").Value = '  -0.16%  '

$ws.Range("D21").Value = '17.52'
$ws.Range("E21").Value = '  +4.90%  '

$ws.Range("D22").Value = '6.431'
$ws.Range("E22").Value = '  +4.27%  '

$ws.Range("D23").Value = '28.312.77'
$ws.Range("E23").Value = '  +1.23%  '

$ws.Range("E24").Value = '  +2.55%  '

$ws.Range("D25").Value = '2.412'
$ws.Range("E25").Value = '  +0.46%  '

$ws.Range("D26").Value = '20.81'
$ws.Range("E26").Value = '  +5.03%  '

$ws.Range("D27").Value = '153.57'
$ws.Range("E27").Value = '  +0.06%  '

$ws.Range("D28").Value = '2.387'
$ws.Range("E28").Value = '  +4.21%  '

$ws.Range("D29").Value = '2.014.01'
$ws.Range("E29").Value = '  +3.18%  '

$ws.Range("D30").Value = '133.55'
$ws.Range("E30").Value = '  +2.04%  '

$ws.Range("D31").Value = '1.264'
$ws.Range("E31").Value = '  -0.84%  '

$ws.Range("D32").Value = '4.035'
$ws.Range("E32").Value = '  +0.09%  '

$ws.Range("D33").Value = '0.09566'
$ws.Range("E33").Value = '  +8.70%  '

$ws.Range("D34").Value = '5.872'
$ws.Range("E34").Value = '  +1.16%  '

$ws.Range("D35").Value = '0.2271'
$ws.Range("E35").Value = '  +8.81%  '

$ws.Range("D36").Value = '0.06401'
$ws.Range("E36").Value = '  +4.69%  '

$ws.Range("D37").Value = '12.13'
$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").Value = '5.292'
$ws.Range("E38").Value = '  +3.41%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.02360'
$ws.Range("E39").Value = '  +3.75%  '

$ws.Range("D40").Value = '0.6661'
$ws.Range("E40").Value = '  +1.87%  '

$ws.Range("D41").Value = '1.237'
$ws.Range("E41").Value = '  +2.59%  '

$ws.Range("D42").Value = '1.491'
$ws.Range("E42").Value = '  -3.19%  '

$ws.Range("D43").Value = '8.328'
$ws.Range("E43").Value = '  +4.71%  '

$ws.Range("D44").Value = '14.30'
$ws.Range("E44").Value = '  +4.16%  '

$ws.Range("D45").Value = '0.9991'
$ws.Range("E45").Value = '  -0.13%  '

$ws.Range("D46").Value = '0.6156'
$ws.Range("E46").Value = '  +2.41%  '

$ws.Range("D47").Value = '3.860'
$ws.Range("E47").Value = '  +0.67%  '

$ws.Range("D48").Value = '130.85'
$ws.Range("E48").Value = '  +3.49%  '

$ws.Range("D49").Value = '2.044'
$ws.Range("E49").Value = '  +2.50%  '

$ws.Range("D50").Value = '1.181'
$ws.Range("E50").Value = '  +1.28%  '

$ws.Range("D51").Value = '0.07155'
$ws.Range("E51").Value = '  +2.51%  '

$ws.Range("D2:D51").ClearFormats()
